$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'249.53"
$ws.Range('D3').Value = "'21.70"
$ws.Range('D4').Value = "'5.448"
$ws.Range('D5').Value = "'0.05690"
$ws.Range('D6').Value = "'3.383"
$ws.Range('D8').Value = "'1.031"
$ws.Range('B9').Value = "'WazirX"
$ws.Range('C9').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D9').Value = "'0.1470"
$ws.Range('E9').Value = "'8WazirXWRX"
$ws.Range('B10').Value = "'MandalaExchangeToken"
$ws.Range('C10').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D10').Value = "'0.07709"
$ws.Range('E10').Value = "'9MandalaExchangeTokenMDX"
$ws.Range('B11').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C11').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D11').Value = "'0.03167"
$ws.Range('E11').Value = "'10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range('B12').Value = "'BitrueCoin"
$ws.Range('C12').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D12').Value = "'0.03026"
$ws.Range('E12').Value = "'11BitrueCoinBTR"
$ws.Range('B13').Value = "'BitMartToken"
$ws.Range('C13').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D13').Value = "'0.09266"
$ws.Range('E13').Value = "'12BitMartTokenBMX"
$ws.Range('B14').Value = "'MCDex"
$ws.Range('C14').Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range('D14').Value = "'3.538"
$ws.Range('E14').Value = "'13MCDexMCB"
$ws.Range('B15').Value = "'BitForexToken"
$ws.Range('C15').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D15').Value = "'0.001653"
$ws.Range('E15').Value = "'14BitForexTokenBF"
$ws.Range('B16').Value = "'CoinExToken"
$ws.Range('C16').Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range('D16').Value = "'0.04713"
$ws.Range('E16').Value = "'15CoinExTokenCET"
$ws.Range('B17').Value = "'One"
$ws.Range('C17').Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range('D17').Value = "'0.0005866"
$ws.Range('E17').Value = "'16OneONE"
$ws.Range('D18').Value = "'0.006364"
$ws.Range('D19').Value = "'0.005033"
$ws.Range('E19').Value = "'18HotbitTokenHTBBestin24h"
$ws.Range('D20').Value = "'0.001043"
$ws.Range('D21').Value = "'0.0001500"
$ws.Range('D23').Value = "'3.771"
$ws.Range('D24').Value = "'6.424"
$ws.Range('D25').Value = "'2.179"
$ws.Range('D26').Value = "'0.3305"
$ws.Range('D40').Value = "'0.04078"
$ws.Range('D41').Value = "'0.006980"
$ws.Range('B42').Value = "'BKEXToken"
$ws.Range('C42').Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range('D42').Value = "'0.1044"
$ws.Range('E42').Value = "'41BKEXTokenBKK"
$ws.Range('B43').Value = "'CEJI"
$ws.Range('C43').Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range('D43').Value = "'0.002881"
$ws.Range('E43').Value = "'42CEJICEJI"
$ws.Range('D44').Value = "'0.007875"
$ws.Range('D45').Value = "'0.00005906"
$ws.Range('D48').Value = "'0.6832"
$ws.Range('D49').Value = "'0.008903"
